# Merge the split "<tl>...</tl> ou <tl>poches de " markup runs around
# "manches" into a single plain-text run, per the refreshed term usage:
#   <tl>manches</tl> ou <tl>poches de   ->   manches ou poches de
# (the surrounding "<tl>" opening tag before "manches" and the
#  "<m>cuir</m></tl>" closing markup after "poches de " are untouched)

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "manches</tl> ou <tl>poches de ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "manches ou poches de ",
    2)

if (-not $found) {
    throw "Target text for the manches/poches merge was not found."
}
